# Updates test results for all Networks
$wb = $excel.ActiveWorkbook

$wsTestData = $wb.Worksheets.Item("TestData")
$wsWindows  = $wb.Worksheets.Item("Windows")

# --- TestData sheet: drop the TC3/Esquire-Network and TC2/Oxygen rows, ---
# --- and change the remaining TC1 row's Network from Oxygen to TeleXitos ---
$wsTestData.Range("E2").Value = "TeleXitos"
$wsTestData.Rows.Item(3).Resize(2, 1).EntireRow.Delete() | Out-Null
$wsTestData.Range("E2").Select() | Out-Null

# --- Windows sheet: TC1 rows become TC2 rows, and new TC1 rows (with ---
# --- updated dates) are appended ---
$wsWindows.Range("A2").Value = "TC2"
$wsWindows.Range("A3").Value = "TC2"

$wsWindows.Range("A4").Value = "TC1"
$wsWindows.Range("C4").Value = 43137
$wsWindows.Range("D4").Value = 43190

$wsWindows.Range("A5").Value = "TC1"
$wsWindows.Range("B5").Value = "W2"
$wsWindows.Range("C5").Value = 43191
$wsWindows.Range("D5").Value = 43251
$wsWindows.Range("E5").Value = 2
$wsWindows.Range("F5").Value = 4

# copy formatting from the row above onto the freshly added row
$wsWindows.Range("A4:F4").Copy() | Out-Null
$wsWindows.Range("A5:F5").PasteSpecial(-4122) | Out-Null

$wsWindows.Range("B6").Select() | Out-Null
